$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3732701332824604
$ws.Range("C2").Value = 0.03891874011020491
$ws.Range("D2").Value = 0.07836328335055498
$ws.Range("E2").Value = 0.1466134964690013
$ws.Range("G2").Value = 0.002454365723921218
$ws.Range("K2").Value = 0.335453974968317
$ws.Range("M2").Value = 0.2437850966233057
$ws.Range("N2").Value = 1.751992835625195
$ws.Range("O2").Value = 3.462580261374399

$ws.Range("B3").Value = 0.3380857983960652
$ws.Range("C3").Value = 0.03448922674714083
$ws.Range("D3").Value = 0.07110752184964042
$ws.Range("E3").Value = 0.1358333083087899
$ws.Range("G3").Value = 0.00245719536263713
$ws.Range("K3").Value = 0.2985031938573002
$ws.Range("M3").Value = 0.2213258291702402
$ws.Range("N3").Value = 1.769541165347839
$ws.Range("O3").Value = 3.458641929997412

$ws.Range("B4").Value = 0.3165956454133436
$ws.Range("C4").Value = 0.03175381070093408
$ws.Range("D4").Value = 0.06668701277178002
$ws.Range("E4").Value = 0.1293086494432316
$ws.Range("G4").Value = 0.002459025663584168
$ws.Range("K4").Value = 0.2758713348815434
$ws.Range("M4").Value = 0.2076408884541507
$ws.Range("N4").Value = 1.78090418464372
$ws.Range("O4").Value = 3.458099837195732

$ws.Range("B5").Value = 0.3078669564978043
$ws.Range("C5").Value = 0.03063517925977521
$ws.Range("D5").Value = 0.06489430411517105
$ws.Range("E5").Value = 0.1266733206934276
$ws.Range("G5").Value = 0.002459794955835244
$ws.Range("K5").Value = 0.2666630073685639
$ws.Range("M5").Value = 0.2020904797138101
$ws.Range("N5").Value = 1.785682543272582
$ws.Range("O5").Value = 3.458350397188951

$ws.Range("B6").Value = 0.3064193069972987
$ws.Range("C6").Value = 0.03044919461501649
$ws.Range("D6").Value = 0.0645971505982601
$ws.Range("E6").Value = 0.1262371406585672
$ws.Range("G6").Value = 0.00245992411349876
$ws.Range("K6").Value = 0.2651348447281521
$ws.Range("M6").Value = 0.2011704244280068
$ws.Range("N6").Value = 1.78648491241427
$ws.Range("O6").Value = 3.458420469786347

$ws.Range("B7").Value = 0.3164778106113033
$ws.Range("C7").Value = 0.03173874032627566
$ws.Range("D7").Value = 0.06666280051857143
$ws.Range("E7").Value = 0.1292730135254843
$ws.Range("G7").Value = 0.002459035943553721
$ws.Range("K7").Value = 0.2757470897802534
$ws.Range("M7").Value = 0.2075659273102133
$ws.Range("N7").Value = 1.780968028800618
$ws.Range("O7").Value = 3.458101307766299

$ws.Range("B8").Value = 0.3611152336411578
$ws.Range("C8").Value = 0.03739470528937261
$ws.Range("D8").Value = 0.07585431894106875
$ws.Range("E8").Value = 0.142876754765112
$ws.Range("G8").Value = 0.002455322146807241
$ws.Range("K8").Value = 0.322701840297924
$ws.Range("M8").Value = 0.236019238317219
$ws.Range("N8").Value = 1.757921350629911
$ws.Range("O8").Value = 3.460832766973823

$ws.Range("B9").Value = 0.4495400192301702
$ws.Range("C9").Value = 0.04836150915005533
$ws.Range("D9").Value = 0.09415426425395879
$ws.Range("E9").Value = 0.1703141902321974
$ws.Range("G9").Value = 0.002448773189791997
$ws.Range("K9").Value = 0.4152190495609034
$ws.Range("M9").Value = 0.2926588579753115
$ws.Range("N9").Value = 1.717397993552328
$ws.Range("O9").Value = 3.481092688501832

$ws.Range("B10").Value = 0.5150462394921362
$ws.Range("C10").Value = 0.05634348765168795
$ws.Range("D10").Value = 0.1077701968604714
$ws.Range("E10").Value = 0.190953906422294
$ws.Range("G10").Value = 0.002444404434575298
$ws.Range("K10").Value = 0.4834582963297009
$ws.Range("M10").Value = 0.3348018616879713
$ws.Range("N10").Value = 1.690476914967569
$ws.Range("O10").Value = 3.505095318791064

$ws.Range("B11").Value = 0.5449639588094897
$ws.Range("C11").Value = 0.05995851081620174
$ws.Range("D11").Value = 0.1140022959725115
$ws.Range("E11").Value = 0.2004519398462463
$ws.Range("G11").Value = 0.002442512143381959
$ws.Range("K11").Value = 0.5145603013347113
$ws.Range("M11").Value = 0.3540925825110648
$ws.Range("N11").Value = 1.678849816047773
$ws.Range("O11").Value = 3.518001996821681

$ws.Range("B12").Value = 0.5563099082238239
$ws.Range("C12").Value = 0.06132511894271886
$ws.Range("D12").Value = 0.1163677412964574
$ws.Range("E12").Value = 0.2040645351704882
$ws.Range("G12").Value = 0.002441809180234773
$ws.Range("K12").Value = 0.5263462557820446
$ws.Range("M12").Value = 0.3614148866847557
$ws.Range("N12").Value = 1.674536119236336
$ws.Range("O12").Value = 3.523175733068882

$ws.Range("B13").Value = 0.5538656144162815
$ws.Range("C13").Value = 0.0610308990370072
$ws.Range("D13").Value = 0.1158580564013505
$ws.Range("E13").Value = 0.2032857874220895
$ws.Range("G13").Value = 0.002441959971659774
$ws.Range("K13").Value = 0.523807574829334
$ws.Range("M13").Value = 0.3598371249332644
$ws.Range("N13").Value = 1.675461179785067
$ws.Range("O13").Value = 3.52204873830479

$ws.Range("B14").Value = 0.5458970627380495
$ws.Range("C14").Value = 0.06007098914895437
$ws.Range("D14").Value = 0.1141967926650835
$ws.Range("E14").Value = 0.2007488303161793
$ws.Range("G14").Value = 0.002442454037741459
$ws.Range("K14").Value = 0.5155297741240474
$ws.Range("M14").Value = 0.3546946449362096
$ws.Range("N14").Value = 1.678493136325272
$ws.Range("O14").Value = 3.518421903777522

$ws.Range("B15").Value = 0.5410182683533264
$ws.Range("C15").Value = 0.0594827135751359
$ws.Range("D15").Value = 0.1131799349057019
$ws.Range("E15").Value = 0.199196947922438
$ws.Range("G15").Value = 0.002442758437553748
$ws.Range("K15").Value = 0.5104604549596559
$ws.Range("M15").Value = 0.3515469874953965
$ws.Range("N15").Value = 1.680361922829022
$ws.Range("O15").Value = 3.5162376545764

$ws.Range("B16").Value = 0.5130934074793458
$ws.Range("C16").Value = 0.05610691340135077
$ws.Range("D16").Value = 0.1073636820494954
$ws.Range("E16").Value = 0.1903353993311967
$ws.Range("G16").Value = 0.002444530008083312
$ws.Range("K16").Value = 0.4814268835825999
$ws.Range("M16").Value = 0.3335435858017632
$ws.Range("N16").Value = 1.691249251377961
$ws.Range("O16").Value = 3.504291871677282

$ws.Range("B17").Value = 0.4959925770815232
$ws.Range("C17").Value = 0.05403185282278855
$ws.Range("D17").Value = 0.103805368615113
$ws.Range("E17").Value = 0.1849271898260625
$ws.Range("G17").Value = 0.00244564111388796
$ws.Range("K17").Value = 0.4636308050601201
$ws.Range("M17").Value = 0.322529814587547
$ws.Range("N17").Value = 1.698087039218429
$ws.Range("O17").Value = 3.497472931848307

$ws.Range("B18").Value = 0.486167813767338
$ws.Range("C18").Value = 0.05283682622768993
$ws.Range("D18").Value = 0.1017623110967634
$ws.Range("E18").Value = 0.1818267691371815
$ws.Range("G18").Value = 0.002446289145531025
$ws.Range("K18").Value = 0.4534006003723903
$ws.Range("M18").Value = 0.3162062401805201
$ws.Range("N18").Value = 1.702078248322522
$ws.Range("O18").Value = 3.493737930178355

$ws.Range("B19").Value = 0.4828432473217674
$ws.Range("C19").Value = 0.05243195257726541
$ws.Range("D19").Value = 0.1010711835564138
$ws.Range("E19").Value = 0.180778771318252
$ws.Range("G19").Value = 0.002446510097713573
$ws.Range("K19").Value = 0.449937803997102
$ws.Range("M19").Value = 0.3140671160169433
$ws.Range("N19").Value = 1.703439610471303
$ws.Range("O19").Value = 3.492505439888333

$ws.Range("B20").Value = 0.4978118345303244
$ws.Range("C20").Value = 0.05425290263448801
$ws.Range("D20").Value = 0.1041837855780017
$ws.Range("E20").Value = 0.1855018416413117
$ws.Range("G20").Value = 0.0024455219084971
$ws.Range("K20").Value = 0.4655246474098931
$ws.Range("M20").Value = 0.3237010835652825
$ws.Range("N20").Value = 1.697353110629493
$ws.Range("O20").Value = 3.498179455324532

$ws.Range("B21").Value = 0.5482371670800887
$ws.Range("C21").Value = 0.0603530011408111
$ws.Range("D21").Value = 0.1146845972470061
$ws.Range("E21").Value = 0.2014935628349193
$ws.Range("G21").Value = 0.002442308549949962
$ws.Range("K21").Value = 0.5179609417700988
$ws.Range("M21").Value = 0.3562046447610001
$ws.Range("N21").Value = 1.677600154145132
$ws.Range("O21").Value = 3.519479420689976

$ws.Range("B22").Value = 0.5812905796811378
$ws.Range("C22").Value = 0.06432622091030282
$ws.Range("D22").Value = 0.1215794529504421
$ws.Range("E22").Value = 0.2120378409926573
$ws.Range("G22").Value = 0.002440287713619623
$ws.Range("K22").Value = 0.5522793438718736
$ws.Range("M22").Value = 0.3775486511713808
$ws.Range("N22").Value = 1.665210650763548
$ws.Range("O22").Value = 3.535068903440589

$ws.Range("B23").Value = 0.5636405244632385
$ws.Range("C23").Value = 0.06220688462860835
$ws.Range("D23").Value = 0.1178966138169102
$ws.Range("E23").Value = 0.206401596170096
$ws.Range("G23").Value = 0.002441359039338815
$ws.Range("K23").Value = 0.5339586422984155
$ws.Range("M23").Value = 0.3661476673412096
$ws.Range("N23").Value = 1.671775510537778
$ws.Range("O23").Value = 3.526595675282124

$ws.Range("B24").Value = 0.4969893272480022
$ws.Range("C24").Value = 0.05415297238778294
$ws.Range("D24").Value = 0.1040126949712885
$ws.Range("E24").Value = 0.185242014068983
$ws.Range("G24").Value = 0.002445575772563435
$ws.Range("K24").Value = 0.4646684381566502
$ws.Range("M24").Value = 0.3231715267417457
$ws.Range("N24").Value = 1.697684732620225
$ws.Range("O24").Value = 3.497859458885443

$ws.Range("B25").Value = 0.4255236620645917
$ws.Range("C25").Value = 0.04540796978729134
$ws.Range("D25").Value = 0.08917384139320461
$ws.Range("E25").Value = 0.1628083795135993
$ws.Range("G25").Value = 0.002450466773947598
$ws.Range("K25").Value = 0.3901437730634711
$ws.Range("M25").Value = 0.2772444527322122
$ws.Range("N25").Value = 1.727860180713751
$ws.Range("O25").Value = 3.474013014783594
